$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.076.08"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "1.818.47"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "310.44"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "0.4969"
$ws.Range("E7").Value = "  -3.84%  "
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "0.09874"
$ws.Range("E9").Value = "  +25.30%  "
$ws.Range("D10").Value = "1.109"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "40.75"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "6.445"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").Value = "20.56"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "1.814.08"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "7.280"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "0.00001140"
$ws.Range("E17").Value = "  +5.32%  "
$ws.Range("D18").Value = "92.34"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "0.06640"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "5.978"
$ws.Range("D23").Value = "28.130.96"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").Value = "11.26"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").Value = "2.239"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "159.33"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "20.77"
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").Value = "2.024.90"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "2.410"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").Value = "126.60"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("E31").Value = "  -3.06%  "
$ws.Range("D32").Value = "1.036"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").Value = "5.573"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Value = "3.612"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("E35").Value = "  -6.32%  "
$ws.Range("D36").Value = "0.02344"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").Value = "8.881"
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").Value = "0.2141"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "4.963"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("D40").Value = "11.37"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("D41").Value = "0.6209"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "1.182"
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").Value = "13.21"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "0.5898"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("D46").Value = "3.697"
$ws.Range("E47").Value = "  -3.81%  "
$ws.Range("D48").Value = "124.81"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").Value = "1.940"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("D51").Value = "0.06772"
$ws.Range("E51").Value = "  -1.38%  "
